$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns
$ws.Range("H1").Value = "lower_case"
$ws.Range("I1").Value = "convert_ascii"

# New data columns (H2:I10)
$values = @(
    @(1, 0),
    @(1, 0),
    @(1, 1),
    @(1, 1),
    @(1, 1),
    @(1, 1),
    @(1, 1),
    @(1, 1),
    @(0, 0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i][0]
    $ws.Cells.Item($row, 9).Value = $values[$i][1]
}

# Drop the now-unused placeholder header cells beyond the two new columns
$ws.Range("J1:L1").Clear()

# Remove the stray formatted row left over below the table (was row 26, H26)
$ws.Rows(26).Delete()

# Update selection to match target workbook
$ws.Range("H3").Select()
